# Insert a new weekly data row before the existing row 153 (shifting the
# remaining "Betarraga" records down by one row, through the former row 246
# which becomes row 247) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(153).Insert()

$ws.Cells.Item(153, 1).Value = 4
$ws.Cells.Item(153, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(153, 3).Value = "Los Lagos"
$ws.Cells.Item(153, 4).Value = 44596
$ws.Cells.Item(153, 5).Value = 10
$ws.Cells.Item(153, 6).Value = 100114014
$ws.Cells.Item(153, 7).Value = "Betarraga"
$ws.Cells.Item(153, 8).Value = "Sin especificar"
$ws.Cells.Item(153, 9).Value = "Primera"
$ws.Cells.Item(153, 10).Value = 1000
$ws.Cells.Item(153, 11).Value = 800
$ws.Cells.Item(153, 12).Value = 900
$ws.Cells.Item(153, 13).Value = 850
$ws.Cells.Item(153, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(153, 15).Value = "Región del Maule"
$ws.Cells.Item(153, 16).Value = 170
$ws.Cells.Item(153, 17).Value = 5
$ws.Cells.Item(153, 18).Value = "Hortaliza"
